# Applies the "imports injection, same message type merge" edit to the
# Reward sheet of the workbook (the commit's meaningful, semantic change).
#
# Context (from the XML diff):
#   - Row 2 of the Reward sheet carries the Go/Proto-ish type annotation for
#     each column. The "InputItem*" columns (C, H) used the generic
#     "[Item]int32" import and now use the more specific "[.Item]int32"
#     import; the "OutputItem1Id" column (N) now uses its own
#     "[OutputItem]int32" import.
#   - Two reward rows that previously duplicated message type "1" (rows 5/6)
#     and an inconsistent row 7 are renumbered 2/3/4 - i.e. same message
#     type rows get merged/renumbered sequentially.
#   - A new reward entry "奖励4" replaces the stale "奖励3" label on row 7's
#     Desc column, matching the renumbered row.
#   - The active selection on the Reward sheet moves from H14 (now out of
#     range) to H9.
#
# Order of writes matters here: the workbook's shared-string table is
# append-only/deduplicating, so writing B7 ("奖励4") before C2/H2/N2
# ("[.Item]int32" / "[OutputItem]int32") reproduces the exact shared-string
# order seen in the target file.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Reward")

# New reward row label (adds shared string "奖励4").
$ws.Cells.Item(7, 2).Value = "奖励4"

# Row 2 type-annotation updates (adds "[.Item]int32" then "[OutputItem]int32").
$ws.Cells.Item(2, 3).Value = "[.Item]int32"   # C2 InputItem1Id
$ws.Cells.Item(2, 8).Value = "[.Item]int32"   # H2 InputItem2Id
$ws.Cells.Item(2, 14).Value = "[OutputItem]int32"  # N2 OutputItem1Id

# Renumber the merged-message-type rows.
$ws.Cells.Item(5, 1).Value = 2   # A5
$ws.Cells.Item(6, 1).Value = 3   # A6
$ws.Cells.Item(7, 1).Value = 4   # A7

# Move the sheet's active selection to match the resaved file.
[void]$ws.Range("H9").Select()
